$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 15
$ws.Range("H15").Value = 1459.3334
$ws.Range("I15").Value = 1459.3334
$ws.Range("K15").Value = 4378.0002
$ws.Range("M15").Value = -4209.0002

# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

# Row 80
$ws.Range("H80").Value = 741.7778
$ws.Range("I80").Value = 573.1667
$ws.Range("J80").Value = 1079
$ws.Range("K80").Value = 1719.5001
$ws.Range("L80").Value = 3237
$ws.Range("M80").Value = -721.5001
$ws.Range("N80").Value = -5233

# Row 83
$ws.Range("H83").Value = 741.7778
$ws.Range("I83").Value = 573.1667
$ws.Range("J83").Value = 1079
$ws.Range("K83").Value = 5158.5003
$ws.Range("L83").Value = 9711
$ws.Range("M83").Value = -166.5002999999997
$ws.Range("N83").Value = -19695

# Row 112
$ws.Range("H112").Value = 1110
$ws.Range("J112").Value = 1110
$ws.Range("L112").Value = 3330
$ws.Range("N112").Value = -5546

# Row 124
$ws.Range("H124").Value = 134373.75
$ws.Range("J124").Value = 134373.75
$ws.Range("L124").Value = 134373.75
$ws.Range("N124").Value = -144193.75

# Row 125
$ws.Range("H125").Value = 6497.5
$ws.Range("I125").Value = 6997.5
$ws.Range("J125").Value = 5997.5
$ws.Range("K125").Value = 62977.5
$ws.Range("L125").Value = 53977.5
$ws.Range("M125").Value = -60517.5
$ws.Range("N125").Value = -58897.5

# Row 129
$ws.Range("H129").Value = 2328
$ws.Range("I129").Value = 2990
$ws.Range("J129").Value = 1997
$ws.Range("K129").Value = 8970
$ws.Range("L129").Value = 5991
$ws.Range("M129").Value = -3970
$ws.Range("N129").Value = -15991

# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 134
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

# Row 135
$ws.Range("H135").Value = 2079.1428
$ws.Range("I135").Value = 2079.1428
$ws.Range("K135").Value = 18712.2852
$ws.Range("M135").Value = -16177.2852

$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 719.05884
$ws.Range("I2").Value = 655.8889
$ws.Range("J2").Value = 790.125
$ws.Range("K2").Value = 655.8889
$ws.Range("L2").Value = 790.125
$ws.Range("M2").Value = -542.8889
$ws.Range("N2").Value = -1016.125

# Row 32
$ws.Range("H32").Value = 3384.72
$ws.Range("I32").Value = 2157.348
$ws.Range("K32").Value = 2157.348
$ws.Range("M32").Value = -1870.348

# Row 116
$ws.Range("H116").Value = 719.05884
$ws.Range("I116").Value = 655.8889
$ws.Range("J116").Value = 790.125
$ws.Range("K116").Value = 655.8889
$ws.Range("L116").Value = 790.125
$ws.Range("M116").Value = 1638.1111
$ws.Range("N116").Value = -5378.125

$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 719.05884
$ws.Range("I3").Value = 655.8889
$ws.Range("J3").Value = 790.125
$ws.Range("K3").Value = 655.8889
$ws.Range("L3").Value = 790.125
$ws.Range("M3").Value = -541.8889
$ws.Range("N3").Value = -1018.125

# Row 57
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# Row 58
$ws.Range("H58").Value = 79999
$ws.Range("J58").Value = 79999
$ws.Range("L58").Value = 79999
$ws.Range("N58").Value = -80587

# Row 74
$ws.Range("H74").Value = 24995
$ws.Range("J74").Value = 24995
$ws.Range("L74").Value = 24995
$ws.Range("N74").Value = -26867

# Row 77
$ws.Range("H77").Value = 24995
$ws.Range("J77").Value = 24995
$ws.Range("L77").Value = 74985
$ws.Range("N77").Value = -84345

# Row 136
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# Row 137
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

$ws = $wb.Worksheets.Item("CRP")

# Row 58
$ws.Range("H58").Value = 2264.5
$ws.Range("I58").Value = 2279
$ws.Range("K58").Value = 2279
$ws.Range("M58").Value = -2076

# Row 86
$ws.Range("H86").Value = 60599.6
$ws.Range("I86").Value = 12749.5
$ws.Range("K86").Value = 12749.5
$ws.Range("M86").Value = -11626.5

# Row 89
$ws.Range("H89").Value = 60599.6
$ws.Range("I89").Value = 12749.5
$ws.Range("K89").Value = 63747.5
$ws.Range("M89").Value = -58131.5

# Row 99
$ws.Range("H99").Value = 1433.3334
$ws.Range("I99").Value = 1150
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1150
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 348
$ws.Range("N99").Value = -4996

# Row 100
$ws.Range("H100").Value = 20000
$ws.Range("J100").Value = 20000
$ws.Range("L100").Value = 20000
$ws.Range("N100").Value = -22164

# Row 126
$ws.Range("H126").Value = 1433.3334
$ws.Range("I126").Value = 1150
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 3450
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -980
$ws.Range("N126").Value = -10940

# Row 129
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# Row 131
$ws.Range("H131").Value = 92197.5
$ws.Range("J131").Value = 92197.5
$ws.Range("L131").Value = 92197.5
$ws.Range("N131").Value = -102277.5

# Row 132
$ws.Range("H132").Value = 4974.125
$ws.Range("I132").Value = 3950
$ws.Range("K132").Value = 11850
$ws.Range("M132").Value = -9320

# Row 136
$ws.Range("H136").Value = 2264.5
$ws.Range("I136").Value = 2279
$ws.Range("K136").Value = 6837
$ws.Range("M136").Value = -4287

# Row 139
$ws.Range("H139").Value = 150000
$ws.Range("J139").Value = 150000
$ws.Range("L139").Value = 150000
$ws.Range("N139").Value = -160280

$ws = $wb.Worksheets.Item("CUL")

# Row 116
$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

# Row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# Row 137
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")

# Row 32
$ws.Range("H32").Value = 45000
$ws.Range("J32").Value = 45000
$ws.Range("L32").Value = 45000
$ws.Range("N32").Value = -45592

# Row 80
$ws.Range("H80").Value = 4764.6665
$ws.Range("I80").Value = 3269.75
$ws.Range("J80").Value = 5960.6
$ws.Range("K80").Value = 3269.75
$ws.Range("L80").Value = 5960.6
$ws.Range("M80").Value = -2271.75
$ws.Range("N80").Value = -7956.6

# Row 83
$ws.Range("H83").Value = 4764.6665
$ws.Range("I83").Value = 3269.75
$ws.Range("J83").Value = 5960.6
$ws.Range("K83").Value = 16348.75
$ws.Range("L83").Value = 29803
$ws.Range("M83").Value = -11356.75
$ws.Range("N83").Value = -39787

# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 132
$ws.Range("H132").Value = 4642.6
$ws.Range("I132").Value = 4303.5
$ws.Range("K132").Value = 12910.5
$ws.Range("M132").Value = -10380.5

$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 5333
$ws.Range("I7").Value = 4999.5
$ws.Range("K7").Value = 4999.5
$ws.Range("M7").Value = -4887.5

# Row 122
$ws.Range("H122").Value = 1999.5
$ws.Range("I122").Value = 1999.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5998.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3548.5
$ws.Range("N122").ClearContents()

# Row 123
$ws.Range("H123").Value = 66650
$ws.Range("J123").Value = 66650
$ws.Range("L123").Value = 66650
$ws.Range("N123").Value = -76450

# Row 126
$ws.Range("H126").Value = 5333
$ws.Range("I126").Value = 4999.5
$ws.Range("K126").Value = 14998.5
$ws.Range("M126").Value = -12528.5
